# Update column F (dSF) values for specific rows as per the re-pulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 9
    3  = -1
    4  = -11
    12 = -4
    13 = -5
    14 = -9
    18 = 6
    21 = -4
    25 = -5
    28 = -4
    31 = -6
    33 = 0
    35 = 0
    38 = 4
    41 = -5
    44 = -1
    45 = 10
    51 = 8
    52 = -6
    53 = 5
    56 = -3
    57 = 9
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
